$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix A26/B26: "IDOCS201" -> "2IDOCS201"
$ws.Range("A26").Value = "2IDOCS201"
$ws.Range("B26").Value = "2IDOCS201"

# Update column E (Annee) from numeric 1/2/3 to text "1A"/"2A"/"3A"
for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 5).Value = "1A"
}
for ($r = 17; $r -le 27; $r++) {
    $ws.Cells.Item($r, 5).Value = "2A"
}
for ($r = 28; $r -le 36; $r++) {
    $ws.Cells.Item($r, 5).Value = "3A"
}

# Update view: scroll and selection
$ws.Range("A26").Select()
